$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($ws, $addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextCell $ws 'D2' '305.57'
Set-TextCell $ws 'E2' '1.29%'
Set-TextCell $ws 'G2' '6'
Set-TextCell $ws 'D3' '36.11'
Set-TextCell $ws 'E3' '-4.13%'
Set-TextCell $ws 'G3' '6'
Set-TextCell $ws 'D4' '5.083'
Set-TextCell $ws 'E4' '1.50%'
Set-TextCell $ws 'G4' '6'
Set-TextCell $ws 'D5' '0.07847'
Set-TextCell $ws 'E5' '-0.13%'
Set-TextCell $ws 'G5' '6'
Set-TextCell $ws 'D6' '2.172'
Set-TextCell $ws 'E6' '-1.65%'
Set-TextCell $ws 'G6' '6'
Set-TextCell $ws 'D7' '7.927'
Set-TextCell $ws 'E7' '-1.16%'
Set-TextCell $ws 'G7' '6'
Set-TextCell $ws 'B8' 'GateToken'
Set-TextCell $ws 'C8' 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
Set-TextCell $ws 'D8' '4.093'
Set-TextCell $ws 'E8' '2.10%'
Set-TextCell $ws 'G8' '6'
Set-TextCell $ws 'B9' 'MXToken'
Set-TextCell $ws 'C9' 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextCell $ws 'D9' '0.9189'
Set-TextCell $ws 'E9' '1.07%'
Set-TextCell $ws 'G9' '6'
Set-TextCell $ws 'B10' 'LiechtensteinCryptoassetsExchange'
Set-TextCell $ws 'C10' 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
Set-TextCell $ws 'D10' '0.09638'
Set-TextCell $ws 'E10' '4.58%'
Set-TextCell $ws 'G10' '6'
Set-TextCell $ws 'B11' 'WazirX'
Set-TextCell $ws 'C11' 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
Set-TextCell $ws 'D11' '0.1863'
Set-TextCell $ws 'E11' '0.47%'
Set-TextCell $ws 'G11' '6'
Set-TextCell $ws 'B12' 'MandalaExchangeToken'
Set-TextCell $ws 'C12' 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
Set-TextCell $ws 'D12' '0.08682'
Set-TextCell $ws 'E12' '2.87%'
Set-TextCell $ws 'G12' '6'
Set-TextCell $ws 'B13' 'BitrueCoin'
Set-TextCell $ws 'C13' 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
Set-TextCell $ws 'D13' '0.03486'
Set-TextCell $ws 'E13' '-0.88%'
Set-TextCell $ws 'G13' '6'
Set-TextCell $ws 'B14' 'BitMartToken'
Set-TextCell $ws 'C14' 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
Set-TextCell $ws 'D14' '0.09918'
Set-TextCell $ws 'E14' '-0.18%'
Set-TextCell $ws 'G14' '6'
Set-TextCell $ws 'B15' 'BitForexToken'
Set-TextCell $ws 'C15' 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
Set-TextCell $ws 'D15' '0.001430'
Set-TextCell $ws 'E15' '-3.09%'
Set-TextCell $ws 'G15' '6'
Set-TextCell $ws 'B16' 'TigerCash'
Set-TextCell $ws 'C16' 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
Set-TextCell $ws 'D16' '0.005685'
Set-TextCell $ws 'E16' '0.68%'
Set-TextCell $ws 'G16' '6'
Set-TextCell $ws 'B17' 'LEO'
Set-TextCell $ws 'C17' 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
Set-TextCell $ws 'D17' '3.463'
Set-TextCell $ws 'E17' '-0.20%'
Set-TextCell $ws 'G17' '6'
Set-TextCell $ws 'E18' '19.24%'
Set-TextCell $ws 'G18' '6'
Set-TextCell $ws 'D19' '0.3425'
Set-TextCell $ws 'E19' '-1.12%'
Set-TextCell $ws 'G19' '6'
Set-TextCell $ws 'D20' '0.1339'
Set-TextCell $ws 'E20' '2.43%'
Set-TextCell $ws 'G20' '6'
Set-TextCell $ws 'D21' '4.851'
Set-TextCell $ws 'E21' '5.83%'
Set-TextCell $ws 'G21' '6'
Set-TextCell $ws 'D22' '0.2200'
Set-TextCell $ws 'E22' '-1.30%'
Set-TextCell $ws 'G22' '6'
Set-TextCell $ws 'D23' '0.04539'
Set-TextCell $ws 'E23' '-2.34%'
Set-TextCell $ws 'G23' '6'
Set-TextCell $ws 'D24' '0.005093'
Set-TextCell $ws 'E24' '14.58%'
Set-TextCell $ws 'G24' '6'
Set-TextCell $ws 'D25' '0.001231'
Set-TextCell $ws 'E25' '0.11%'
Set-TextCell $ws 'G25' '6'
Set-TextCell $ws 'D26' '0.0001400'
Set-TextCell $ws 'E26' '7.74%'
Set-TextCell $ws 'G26' '6'
Set-TextCell $ws 'D27' '0.0004749'
Set-TextCell $ws 'E27' '0.09%'
Set-TextCell $ws 'G27' '6'
Set-TextCell $ws 'G28' '6'
Set-TextCell $ws 'G29' '6'
Set-TextCell $ws 'G30' '6'
Set-TextCell $ws 'G31' '6'
Set-TextCell $ws 'G32' '6'
Set-TextCell $ws 'G33' '6'
Set-TextCell $ws 'G34' '6'
Set-TextCell $ws 'G35' '6'
Set-TextCell $ws 'G36' '6'
Set-TextCell $ws 'G37' '6'
Set-TextCell $ws 'G38' '6'
Set-TextCell $ws 'D39' '0.01820'
Set-TextCell $ws 'E39' '3.66%'
Set-TextCell $ws 'G39' '6'
Set-TextCell $ws 'D40' '0.04766'
Set-TextCell $ws 'E40' '1.12%'
Set-TextCell $ws 'G40' '6'
Set-TextCell $ws 'D41' '0.007717'
Set-TextCell $ws 'E41' '-2.61%'
Set-TextCell $ws 'G41' '6'
Set-TextCell $ws 'D42' '0.1396'
Set-TextCell $ws 'E42' '0.37%'
Set-TextCell $ws 'G42' '6'
Set-TextCell $ws 'D43' '0.007746'
Set-TextCell $ws 'E43' '1.07%'
Set-TextCell $ws 'G43' '6'
Set-TextCell $ws 'D44' '0.002215'
Set-TextCell $ws 'E44' '1.18%'
Set-TextCell $ws 'G44' '6'
Set-TextCell $ws 'D45' '0.01119'
Set-TextCell $ws 'E45' '7.96%'
Set-TextCell $ws 'G45' '6'
Set-TextCell $ws 'D46' '0.00006380'
Set-TextCell $ws 'E46' '6.39%'
Set-TextCell $ws 'G46' '6'
Set-TextCell $ws 'D47' '0.00000000750'
Set-TextCell $ws 'E47' '0.10%'
Set-TextCell $ws 'G47' '6'
Set-TextCell $ws 'D48' '0.0005800'
Set-TextCell $ws 'E48' '-0.01%'
Set-TextCell $ws 'G48' '6'
Set-TextCell $ws 'D49' '24.55'
Set-TextCell $ws 'E49' '183.17%'
Set-TextCell $ws 'G49' '6'
Set-TextCell $ws 'D50' '0.002000'
Set-TextCell $ws 'E50' '-25.90%'
Set-TextCell $ws 'G50' '6'
Set-TextCell $ws 'D51' '0.00002100'
Set-TextCell $ws 'E51' '0.10%'
Set-TextCell $ws 'G51' '6'

Write-Output "Applied 147 cell updates"